$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '27.568.80'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.50%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.623.13'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -1.40%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '211.65'
$ws.Range('D5').Style = "Normal"
$ws.Range('E6').Value = '  -0.87%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '23.21'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.262'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +1.61%  '
$ws.Range('E10').Value = '  -0.18%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0890'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -0.24%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.853.43'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -1.38%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.611.41'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -2.36%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.04'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.01%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.549'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -2.28%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '65.30'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.80%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '27.532.52'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.58%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '231.59'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.12%  '
$ws.Range('E19').Value = '  -0.85%  '
$ws.Range('E20').Value = '  -1.24%  '
$ws.Range('E21').Value = '  -0.07%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '10.41'
$ws.Range('D22').Style = "Normal"
$ws.Range('E23').Value = '  +0.81%  '
$ws.Range('E24').Value = '  +5.73%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '149.64'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.34%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '6.87'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.91%  '
$ws.Range('E27').Value = '  -0.94%  '
$ws.Range('B28').Value = 'BinanceUSD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.00'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.54'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.72%  '
$ws.Range('E30').Value = '  -0.99%  '
$ws.Range('E31').Value = '  -0.53%  '
$ws.Range('E32').Value = '  -0.89%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.465.18'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +1.72%  '
$ws.Range('E34').Value = '  -2.42%  '
$ws.Range('E35').Value = '  -2.69%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.34'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +0.08%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.943'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +6.74%  '
$ws.Range('E38').Value = '  +0.59%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.873'
$ws.Range('D39').Style = "Normal"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.554'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -2.97%  '
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('E42').Value = '  -1.86%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '67.55'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.09%  '
$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.20'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -2.08%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '5.29'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -5.63%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.76'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +1.56%  '
$ws.Range('B47').Value = 'RocketPoolETH'
$ws.Range('C47').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.763.79'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -1.37%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '87.41'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +2.11%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0105'
$ws.Range('E49').Value = '  -1.99%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.100'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.54%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '7.68'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.86%  '
